$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dynamics")

$ws.Range("B2").Value = "k:/github/digitalmodel/tests/modules/orcaflex/orcaflex_post_process/orcaflex_test1.sim"
$ws.Range("D2").Value = "k:/github/digitalmodel/tests/modules/orcaflex/orcaflex_post_process/orcaflex_test1.sim"
$ws.Range("B3").Value = "k:/github/digitalmodel/tests/modules/orcaflex/orcaflex_post_process/orcaflex_test2.sim"
$ws.Range("D3").Value = "k:/github/digitalmodel/tests/modules/orcaflex/orcaflex_post_process/orcaflex_test2.sim"
